# Apply updated cryptocurrency price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.805.65"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "2.241.40"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'114.60"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").Value = "'269.30"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("D8").Value = "'1.01"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'46.27"
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("D11").Value = "'0.0929"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "'9.18"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "2.576.69"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "2.239.15"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "43.031.82"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "'72.10"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'2.35"
$ws.Range("E22").Value = "  -5.03%  "
$ws.Range("D23").Value = "'232.69"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").Value = "'9.37"
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").Value = "'12.25"
$ws.Range("E26").Value = "  +7.02%  "
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "'40.68"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").Value = "'173.84"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'21.18"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("E35").Value = "  +10.12%  "
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").Value = "'4.68"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("D41").Value = "'71.45"
$ws.Range("E41").Value = "  -6.35%  "
$ws.Range("D42").Value = "'13.24"
$ws.Range("E42").Value = "  -6.98%  "
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'5.68"
$ws.Range("E45").Value = "  -8.71%  "
$ws.Range("D46").Value = "'1.33"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("D47").Value = "'1.25"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "'8.47"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").Value = "'0.646"
$ws.Range("E50").Value = "  +8.23%  "
$ws.Range("E51").Value = "  -3.35%  "
